$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.449.14"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "1.676.20"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'216.44"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").Value = "'0.5308"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +4.37%  "
$ws.Range("D9").Value = "'0.06393"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").Value = "'21.56"
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("D12").Value = "1.681.83"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").Value = "'4.501"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "'0.5564"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "0.0₅8341"
$ws.Range("E15").Value = "  +6.02%  "
$ws.Range("D16").Value = "'65.63"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "26.477.00"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'4.768"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("D20").Value = "'195.27"
$ws.Range("E20").Value = "  +6.80%  "
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").Value = "'6.335"
$ws.Range("E22").Value = "  +4.95%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'144.17"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'0.1276"
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").Value = "'7.423"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +5.55%  "
$ws.Range("D28").Value = "'1.429"
$ws.Range("E28").Value = "  +5.81%  "
$ws.Range("D29").Value = "'0.06129"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").Value = "'1.274"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").Value = "'3.624"
$ws.Range("E31").Value = "  +7.87%  "
$ws.Range("D32").Value = "'3.447"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("E33").Value = "  +6.08%  "
$ws.Range("D34").Value = "'1.004"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "'2.776"
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("D37").Value = "'0.5724"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("D39").Value = "'6.041"
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("D40").Value = "1.072.43"
$ws.Range("E40").Value = "  +5.57%  "
$ws.Range("D41").Value = "'0.8590"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "1.825.69"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("D45").Value = "'57.05"
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("D46").Value = "0.0₈104"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'8.145"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'0.05205"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "'1.477"
$ws.Range("E50").Value = "  +7.79%  "
$ws.Range("D51").Value = "'6.024"
$ws.Range("E51").Value = "  +4.40%  "
